$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.903.83"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.673.54"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.531"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.24%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.25%  "
$ws.Range("E11").Value = "  +3.74%  "
$ws.Range("D12").Value = "1.910.11"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").Value = "1.673.44"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "26.951.91"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.24%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  +2.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "1.462.89"
$ws.Range("E33").Value = "  -3.76%  "
$ws.Range("E34").Value = "  +4.68%  "
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.907"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.90%  "
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.972"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.88%  "
$ws.Range("D45").Value = "1.819.18"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  +0.67%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.59%  "
